$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 additions (Analysis_Number, Prey_Kingdom, Prey_Phylum, Prey_Class, Bird_Sample_Size)
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = "Animalia"
$ws.Range("U2").Value = "Arthropoda"
$ws.Range("V2").Value = "Insecta"
$ws.Range("AK2").Value = 57

# Row 3 (prey taxonomy continued)
$ws.Range("T3").Value = "Animalia"
$ws.Range("U3").Value = "Arthropoda"
$ws.Range("V3").Value = "Malacostraca"
$ws.Range("W3").Value = "Decapoda"

# Row 4 (prey taxonomy continued)
$ws.Range("T4").Value = "Animalia"
$ws.Range("U4").Value = "Chordata"

# Rows 5-26: new bird species diet records
# Row 5
$ws.Range("A5").Value = "Pied-billed Grebe"
$ws.Range("B5").Value = "Podilymbus podiceps"
$ws.Range("D5").Value = "Podicipedidae"
$ws.Range("E5").Value = "eBird Clements Checklist v2018"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "NA"
$ws.Range("H5").Value = "NA"
$ws.Range("I5").Value = "NA"
$ws.Range("J5").Value = "NA"
$ws.Range("K5").Value = "Alabama"
# Row 6
$ws.Range("A6").Value = "Common Loon"
$ws.Range("B6").Value = "Gavia immer"
$ws.Range("D6").Value = "Gaviidae"
$ws.Range("E6").Value = "eBird Clements Checklist v2018"
# Row 7
$ws.Range("A7").Value = "European Herring Gull"
$ws.Range("B7").Value = "Larus argentatus"
$ws.Range("D7").Value = "Laridae"
$ws.Range("E7").Value = "eBird Clements Checklist v2018"
# Row 8
$ws.Range("A8").Value = "Ring-billed Gull"
$ws.Range("B8").Value = "Larus delawarensis"
$ws.Range("D8").Value = "Laridae"
$ws.Range("E8").Value = "eBird Clements Checklist v2018"
# Row 9
$ws.Range("A9").Value = "Laughing Gull"
$ws.Range("B9").Value = "Leucophaeus atricilla"
$ws.Range("D9").Value = "Laridae"
$ws.Range("E9").Value = "eBird Clements Checklist v2018"
# Row 10
$ws.Range("A10").Value = "Bonaparte's Gull"
$ws.Range("B10").Value = "Chroicocephalus philadelphia"
$ws.Range("D10").Value = "Laridae"
$ws.Range("E10").Value = "eBird Clements Checklist v2018"
# Row 11
$ws.Range("A11").Value = "Gull-billed Tern"
$ws.Range("B11").Value = "Gelochelidon nilotica"
$ws.Range("D11").Value = "Laridae"
$ws.Range("E11").Value = "eBird Clements Checklist v2018"
# Row 12
$ws.Range("A12").Value = "Caspian Tern"
$ws.Range("B12").Value = "Hydroprogne caspia"
$ws.Range("D12").Value = "Laridae"
$ws.Range("E12").Value = "eBird Clements Checklist v2018"
# Row 13
$ws.Range("A13").Value = "Royal Tern"
$ws.Range("B13").Value = "Thalasseus maximus"
$ws.Range("D13").Value = "Laridae"
$ws.Range("E13").Value = "eBird Clements Checklist v2018"
# Row 14
$ws.Range("A14").Value = "Cabot's Tern"
$ws.Range("B14").Value = "Thalasseus acuflavidus "
$ws.Range("D14").Value = "Laridae"
$ws.Range("E14").Value = "eBird Clements Checklist v2018"
# Row 15
$ws.Range("A15").Value = "Forster's Tern"
$ws.Range("B15").Value = "Sterna forsteri"
$ws.Range("D15").Value = "Laridae"
$ws.Range("E15").Value = "eBird Clements Checklist v2018"
# Row 16
$ws.Range("A16").Value = "Common Tern"
$ws.Range("B16").Value = "Sterna hirundo"
$ws.Range("D16").Value = "Laridae"
$ws.Range("E16").Value = "eBird Clements Checklist v2018"
# Row 17
$ws.Range("A17").Value = "Little Tern"
$ws.Range("B17").Value = "Sterna albifrons"
$ws.Range("D17").Value = "Laridae"
$ws.Range("E17").Value = "eBird Clements Checklist v2018"
# Row 18
$ws.Range("A18").Value = "Black Tern"
$ws.Range("B18").Value = "Chlidonias niger"
$ws.Range("D18").Value = "Laridae"
$ws.Range("E18").Value = "eBird Clements Checklist v2018"
# Row 19
$ws.Range("A19").Value = "Black Skimmer"
$ws.Range("B19").Value = "Rynchops niger"
$ws.Range("D19").Value = "Laridae"
$ws.Range("E19").Value = "eBird Clements Checklist v2018"
# Row 20
$ws.Range("A20").Value = "Sooty Shearwater"
$ws.Range("B20").Value = "Puffinus griseus"
$ws.Range("D20").Value = "Laridae"
$ws.Range("E20").Value = "eBird Clements Checklist v2018"
# Row 21
$ws.Range("A21").Value = "Mallard"
$ws.Range("B21").Value = "Anas platyrhnchos"
$ws.Range("D21").Value = "Laridae"
$ws.Range("E21").Value = "eBird Clements Checklist v2018"
# Row 22
$ws.Range("A22").Value = "American Black Duck"
$ws.Range("B22").Value = "Anas rubripes rubripes"
$ws.Range("D22").Value = "Anatidae"
$ws.Range("E22").Value = "eBird Clements Checklist v2018"
# Row 23
$ws.Range("A23").Value = "Gadwall"
$ws.Range("B23").Value = "Anas strepera"
$ws.Range("D23").Value = "Anatidae"
$ws.Range("E23").Value = "eBird Clements Checklist v2018"
# Row 24
$ws.Range("A24").Value = "Baldpate"
$ws.Range("B24").Value = "Mareca americana"
$ws.Range("D24").Value = "Anatidae"
$ws.Range("E24").Value = "eBird Clements Checklist v2018"
# Row 25
$ws.Range("A25").Value = "Green-winged Teal"
$ws.Range("B25").Value = "Anas carolinensis"
$ws.Range("D25").Value = "Anatidae"
$ws.Range("E25").Value = "eBird Clements Checklist v2018"
# Row 26
$ws.Range("A26").Value = "Blue-winged Teal"
$ws.Range("B26").Value = "Spatula discors"
$ws.Range("D26").Value = "Anatidae"
$ws.Range("E26").Value = "eBird Clements Checklist v2018"

# Update sheet view selection and scroll position
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("V4").Select()
